# Insert two new weekly price rows for "Pimiento" (Zafiro rojo / Zafiro verde)
# at the top of the historical data block, pushing the existing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before row 206 (formatting of row 206 carries to the new rows).
$ws.Rows.Item(206).Resize(2).Insert()

# New row 206: Pimiento, Zafiro rojo, Primera - 2022-07-13
$ws.Cells.Item(206, 1).Value = 7
$ws.Cells.Item(206, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(206, 3).Value = "Ñuble"
$ws.Cells.Item(206, 4).Value = 44755
$ws.Cells.Item(206, 5).Value = 16
$ws.Cells.Item(206, 6).Value = 100112002
$ws.Cells.Item(206, 7).Value = "Pimiento"
$ws.Cells.Item(206, 8).Value = "Zafiro rojo"
$ws.Cells.Item(206, 9).Value = "Primera"
$ws.Cells.Item(206, 10).Value = 40
$ws.Cells.Item(206, 11).Value = 35000
$ws.Cells.Item(206, 12).Value = 35000
$ws.Cells.Item(206, 13).Value = 35000
$ws.Cells.Item(206, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(206, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(206, 16).Value = 2333
$ws.Cells.Item(206, 17).Value = 15
$ws.Cells.Item(206, 18).Value = "Hortaliza"

# New row 207: Pimiento, Zafiro verde, Primera - 2022-07-13
$ws.Cells.Item(207, 1).Value = 7
$ws.Cells.Item(207, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(207, 3).Value = "Ñuble"
$ws.Cells.Item(207, 4).Value = 44755
$ws.Cells.Item(207, 5).Value = 16
$ws.Cells.Item(207, 6).Value = 100112002
$ws.Cells.Item(207, 7).Value = "Pimiento"
$ws.Cells.Item(207, 8).Value = "Zafiro verde"
$ws.Cells.Item(207, 9).Value = "Primera"
$ws.Cells.Item(207, 10).Value = 40
$ws.Cells.Item(207, 11).Value = 26000
$ws.Cells.Item(207, 12).Value = 26000
$ws.Cells.Item(207, 13).Value = 26000
$ws.Cells.Item(207, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(207, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(207, 16).Value = 1733
$ws.Cells.Item(207, 17).Value = 15
$ws.Cells.Item(207, 18).Value = "Hortaliza"
